$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 100
$ws.Range("T2").Value = ""
$ws.Range("W2").Value = 100
$ws.Range("X2").Value = ""
$ws.Range("AA2").Value = 100
$ws.Range("AB2").Value = ""
$ws.Range("AE2").Value = 90
$ws.Range("AF2").Value = ""
$ws.Range("AM2").Value = 90
$ws.Range("AN2").Value = ""
$ws.Range("BC2").Value = 70

$ws.Range("BC3").Select()
